$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add "Wins", "Losses", "Ties" in AD1:AF1 ---
# Copy the formatting of the existing header cell (AC1) onto the new header
# cells so they pick up the same style (bold, bordered, centered) without
# creating brand-new style entries.
$headerSrc = $ws.Range("AC1")
$headerSrc.Copy()
$newHeaders = $ws.Range("AD1:AF1")
$newHeaders.PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2 through 68): season record for each player's team ---
$wins = 62
$losses = 100
$ties = 0

for ($r = 2; $r -le 68; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
